$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5450
$ws1.Range("F4").Value = 11816
$ws1.Range("F5").Value = 290
$ws1.Range("F6").Value = 598
$ws1.Range("F7").Value = 173
$ws1.Range("F8").Value = 296
$ws1.Range("F9").Value = 1062
$ws1.Range("F10").Value = 101

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F4").Value = 5450
$ws4.Range("F7").Value = 11816
$ws4.Range("F8").Value = 290
$ws4.Range("F9").Value = 598
$ws4.Range("F10").Value = 173
$ws4.Range("F13").Value = 296
$ws4.Range("F14").Value = 1062
$ws4.Range("F16").Value = 101
